$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date in column C for rows 2-5 from 45175 (2023-09-06)
# to 45183 (2023-09-14), keeping the existing date serial/format.
$ws.Range("C2:C5").Value = 45183
